# Elimina antiguos EC y agrega nuevos y modifica Antigua BD
#
# Replaces the single worker record on the "Antigua BD" row (row 16) with
# a new worker's data: document number, name and overdue-period code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Old record: 1047423111 / ANDREINA CERA PAYARES / 2508
# New record: 1007360236 / NIXA MARIA RIVERA GOMEZ / 2509
$ws.Range("C16").Value = "1007360236"
$ws.Range("D16").Value = "NIXA MARIA RIVERA GOMEZ"
$ws.Range("E16").Value = "2509"

# The "Periodo Mora" cell is given a centered horizontal alignment to match
# the rest of the data row.
$ws.Range("E16").HorizontalAlignment = -4108

# Column D ("Nombre Trabajador") auto-fits to the new, longer name.
$ws.Columns.Item(4).AutoFit()
